# "Generate Report for Handoff"
#
# Updates the localization-status report with a new handoff event for the
# "728a2081-30b5-4860-9087-31d6e377a9a7" source file: a fresh "Latest
# Handoff" timestamp is recorded on the Overview sheet and on each
# per-locale sheet (zh-cn, de-de) for that file's row.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 7 is the 728a2081... file, column D is
#     "Latest Handoff Date" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D7").Value = "2016-36-20 16:36:07"

# --- zh-cn sheet: row 7 is the 728a2081... file, column E is
#     "Latest Handoff Datetime" ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E7").Value = "2016-03-20 16:36:04"

# --- de-de sheet: row 7 is the 728a2081... file. Column D
#     ("Latest Handoff File") and column E ("Latest Handoff Datetime")
#     exchange their values: the handoff-file name moves into E, and the
#     newly recorded handoff timestamp goes into D. ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D7").Value = "2016-03-20 16:36:07"
$dede.Range("E7").Value = "728a2081-30b5-4860-9087-31d6e377a9a7.747ead6bf64fd99a13a9cb00cf239beadec31735.de-de.xlf"
